$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The second paragraph originally reads (as two runs split by the
# "_GoBack" bookmark):
#   Run A: "In this game, ... left and right arrow keys"
#   <bookmarkStart/><bookmarkEnd name="_GoBack">
#   Run B: ".  The game is designed ... 9 lives."
#
# The edit appends new sentences right after "arrow keys" (merging
# the old Run B content back in, plus new continuation sentences and
# the start of a new sentence), keeps the "_GoBack" bookmark exactly
# where it was (immediately before the word "programmers"), and
# replaces the remainder of the paragraph with a new sentence about
# buggy programmers.
# ------------------------------------------------------------------

$oldRunBText = ".  The game is designed such that the Lightning sprite will point to the cat before moving towards it.  If the Lightning hits the cat, it loses a life.  If it avoids getting hit by the lightning, it scores a point.  The game will start with the cat having 9 lives."

$insertedBeforeBookmark = $oldRunBText + "  " + "Unfortunately, the two "

$newTextAfterBookmark = "programmers you hired to write this software were not as great of a programmer as you are and made some errors.  Your task is to find and fix those bugs.  There is only one bug in the program.  This is a process known as debugging.  The first step in debugging is to recognize the symptom of what the program is doing incorrectly.  Sometimes this can be a bit tricky in that there may be more than one symptom from an underlying cause.  They only way to get good at debugging is to practice it a lot."

# Locate the (empty/collapsed) "_GoBack" bookmark that sits right at the
# boundary between the two runs and insert the new text immediately
# before it -- this naturally extends/merges the preceding run, exactly
# as Word does when you place your cursor there and type, while leaving
# the bookmark anchored at the same logical spot in the text.
$bm = $d.Bookmarks("_GoBack")
$insertPoint = $d.Range($bm.Start, $bm.Start)
$insertPoint.InsertBefore($insertedBeforeBookmark)

# Re-resolve the bookmark (its position shifted because we inserted text
# before it) and the paragraph so we can replace everything from right
# after the bookmark through the end of the paragraph (but not the
# paragraph mark itself) with the new "buggy programmers" sentence.
$bm = $d.Bookmarks("_GoBack")
$para = $d.Paragraphs(2)
$paraEnd = $para.Range.End

$tailRange = $d.Range($bm.End, $paraEnd - 1)
$tailRange.Text = $newTextAfterBookmark

Write-Output $d.Paragraphs(2).Range.Text
